# Applies the inventory-update diff:
#  - updates product metrics for existing rows 2-15
#  - inserts two new product rows (16: "25-954C-QWS0"/New Wipes, 17: "SEIWBW"/Body Wash)
#  - refreshes the trailing Total row (now row 18) with the new aggregate figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Cell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# --- 1) Refresh existing product rows (2-15) with the new monthly figures ---
# Row 2: Classic
Set-Cell 2 3 "Classic"
Set-Cell 2 10 4531
Set-Cell 2 11 12
Set-Cell 2 12 6
Set-Cell 2 13 10
Set-Cell 2 14 16
Set-Cell 2 15 4537
Set-Cell 2 16 4286
Set-Cell 2 17 504.85
Set-Cell 2 18 3974
Set-Cell 2 19 1
Set-Cell 2 20 0
Set-Cell 2 21 558
Set-Cell 2 22 4531

# Row 3: Passion Fruit
Set-Cell 3 3 "Passion Fruit"
Set-Cell 3 10 838
Set-Cell 3 11 99
Set-Cell 3 12 290
Set-Cell 3 13 1
Set-Cell 3 14 6
Set-Cell 3 15 450
Set-Cell 3 16 33101
Set-Cell 3 17 70.73
Set-Cell 3 18 772
Set-Cell 3 19 0
Set-Cell 3 20 0
Set-Cell 3 21 66
Set-Cell 3 22 838

# Row 4: Menthol
Set-Cell 4 3 "Menthol"
Set-Cell 4 10 228
Set-Cell 4 11 0
Set-Cell 4 12 0
Set-Cell 4 13 3
Set-Cell 4 14 2
Set-Cell 4 15 232
Set-Cell 4 16 46507
Set-Cell 4 17 30.69
Set-Cell 4 18 202
Set-Cell 4 19 0
Set-Cell 4 20 0
Set-Cell 4 21 26
Set-Cell 4 22 228

# Row 5: Women
Set-Cell 5 3 "Women"
Set-Cell 5 10 58
Set-Cell 5 11 0
Set-Cell 5 12 0
Set-Cell 5 13 0
Set-Cell 5 14 0
Set-Cell 5 15 58
Set-Cell 5 16 148209
Set-Cell 5 17 7.88
Set-Cell 5 18 53
Set-Cell 5 19 0
Set-Cell 5 20 0
Set-Cell 5 21 5
Set-Cell 5 22 58

# Row 6: Intimate Wipes
Set-Cell 6 3 "Intimate Wipes"
Set-Cell 6 10 258
Set-Cell 6 11 1
Set-Cell 6 12 262
Set-Cell 6 13 0
Set-Cell 6 14 0
Set-Cell 6 15 0
Set-Cell 6 16 23989
Set-Cell 6 17 31.77
Set-Cell 6 18 216
Set-Cell 6 19 0
Set-Cell 6 20 0
Set-Cell 6 21 42
Set-Cell 6 22 258

# Row 7: Refill Pack
Set-Cell 7 3 "Refill Pack"
Set-Cell 7 10 779
Set-Cell 7 11 43
Set-Cell 7 12 749
Set-Cell 7 13 0
Set-Cell 7 14 0
Set-Cell 7 15 0
Set-Cell 7 16 10746
Set-Cell 7 17 91.12
Set-Cell 7 18 631
Set-Cell 7 19 0
Set-Cell 7 20 0
Set-Cell 7 21 148
Set-Cell 7 22 779

# Row 8: Turmeric
Set-Cell 8 3 "Turmeric"
Set-Cell 8 10 318
Set-Cell 8 11 0
Set-Cell 8 12 3
Set-Cell 8 13 1
Set-Cell 8 14 1
Set-Cell 8 15 313
Set-Cell 8 16 157552
Set-Cell 8 17 35.74
Set-Cell 8 18 300
Set-Cell 8 19 0
Set-Cell 8 20 0
Set-Cell 8 21 18
Set-Cell 8 22 318

# Row 9: Shampoo
Set-Cell 9 3 "Shampoo"
Set-Cell 9 10 0
Set-Cell 9 18 0
Set-Cell 9 19 0
Set-Cell 9 20 0
Set-Cell 9 21 0
Set-Cell 9 22 0

# Row 10: Classic +Classic
Set-Cell 10 3 "Classic +Classic"
Set-Cell 10 10 944
Set-Cell 10 11 54
Set-Cell 10 12 828
Set-Cell 10 13 67
Set-Cell 10 14 0
Set-Cell 10 15 0
Set-Cell 10 16 4286
Set-Cell 10 17 112.43
Set-Cell 10 18 866
Set-Cell 10 19 0
Set-Cell 10 20 0
Set-Cell 10 21 78
Set-Cell 10 22 944

# Row 11: Classic+ Passion Fruit
Set-Cell 11 3 "Classic+ Passion Fruit"
Set-Cell 11 10 1051
Set-Cell 11 11 92
Set-Cell 11 12 966
Set-Cell 11 13 0
Set-Cell 11 14 0
Set-Cell 11 15 0
Set-Cell 11 16 4286
Set-Cell 11 17 97.64
Set-Cell 11 18 929
Set-Cell 11 19 0
Set-Cell 11 20 0
Set-Cell 11 21 122
Set-Cell 11 22 1051

# Row 12: Classic + Menthol
Set-Cell 12 3 "Classic + Menthol"
Set-Cell 12 10 507
Set-Cell 12 11 54
Set-Cell 12 12 396
Set-Cell 12 13 58
Set-Cell 12 14 0
Set-Cell 12 15 0
Set-Cell 12 16 0
Set-Cell 12 17 61.91
Set-Cell 12 18 504
Set-Cell 12 19 0
Set-Cell 12 20 0
Set-Cell 12 21 3
Set-Cell 12 22 507

# Row 13: Classic + Wipes
Set-Cell 13 3 "Classic + Wipes"
Set-Cell 13 10 102
Set-Cell 13 11 6
Set-Cell 13 12 96
Set-Cell 13 13 2
Set-Cell 13 14 0
Set-Cell 13 15 0
Set-Cell 13 16 165150
Set-Cell 13 17 19.88
Set-Cell 13 18 99
Set-Cell 13 19 0
Set-Cell 13 20 0
Set-Cell 13 21 3
Set-Cell 13 22 102

# Row 14: Wipes + Passion Fruit
Set-Cell 14 3 "Wipes + Passion Fruit"
Set-Cell 14 10 0
Set-Cell 14 18 0
Set-Cell 14 19 0
Set-Cell 14 20 0
Set-Cell 14 21 0
Set-Cell 14 22 0

# Row 15: Wipes + Menthol
Set-Cell 15 3 "Wipes + Menthol"
Set-Cell 15 10 141
Set-Cell 15 11 18
Set-Cell 15 12 124
Set-Cell 15 13 0
Set-Cell 15 14 0
Set-Cell 15 15 0
Set-Cell 15 16 112992
Set-Cell 15 17 42.53
Set-Cell 15 18 110
Set-Cell 15 19 0
Set-Cell 15 20 0
Set-Cell 15 21 31
Set-Cell 15 22 141

# --- 2) Insert two new product rows before the Total row (old row 16 -> shifts to row 18) ---
$ws.Range("A16:A17").EntireRow.Insert()

# --- 3) Populate the new row 16: SKU 25-954C-QWS0 / "New Wipes" ---
Set-Cell 16 1 15
Set-Cell 16 2 "25-954C-QWS0"
Set-Cell 16 3 "New Wipes"
Set-Cell 16 5 0
Set-Cell 16 10 885
Set-Cell 16 11 154
Set-Cell 16 12 630
Set-Cell 16 13 105
Set-Cell 16 14 0
Set-Cell 16 15 0
Set-Cell 16 16 23989
Set-Cell 16 17 76.12
Set-Cell 16 18 784
Set-Cell 16 19 0
Set-Cell 16 20 0
Set-Cell 16 21 101
Set-Cell 16 22 885

# --- 4) Populate the new row 17: SKU SEIWBW / "Body Wash" ---
Set-Cell 17 1 16
Set-Cell 17 2 "SEIWBW"
Set-Cell 17 3 "Body Wash"
Set-Cell 17 5 0
Set-Cell 17 10 83
Set-Cell 17 11 83
Set-Cell 17 12 0
Set-Cell 17 13 0
Set-Cell 17 14 0
Set-Cell 17 15 0
Set-Cell 17 16 54783
Set-Cell 17 17 6.16
Set-Cell 17 18 47
Set-Cell 17 19 0
Set-Cell 17 20 0
Set-Cell 17 21 36
Set-Cell 17 22 83

# --- 5) Refresh the Total row, now shifted to row 18 ---
Set-Cell 18 3 "Total"
Set-Cell 18 5 0
Set-Cell 18 10 10723
Set-Cell 18 11 616
Set-Cell 18 12 4350
Set-Cell 18 13 247
Set-Cell 18 14 25
Set-Cell 18 15 5590
Set-Cell 18 16 789876
Set-Cell 18 17 1189.45
Set-Cell 18 18 9487
Set-Cell 18 19 1
Set-Cell 18 20 0
Set-Cell 18 21 1237
Set-Cell 18 22 10723

